# Update AMF_Late_PercentAllocation data to reflect the revised simulation
# constraints, and add the (currently blank, formatted) columns M:W that
# will hold the data for the new "adding a nodule to a plant that has/
# doesn't have AMF" figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Updated percent-allocation values (rows 2-11, columns A-K)
# ---------------------------------------------------------------------
$data = @{
    2  = @(100, 3.78, 4.41, 4.5,  4.51, 4.51, 4.51, 4.51, 1.46, 0.09, 0)
    3  = @(90,  3.78, 4.41, 4.51, 4.51, 4.51, 4.51, 4.51, 1.46, 0.09, 0.01)
    4  = @(80,  3.78, 4.41, 4.51, 4.51, 4.51, 4.51, 4.51, 1.46, 0.08, 0.02)
    5  = @(70,  3.78, 4.41, 4.51, 4.51, 4.51, 4.51, 4.51, 1.47, 0.07, 0.03)
    6  = @(60,  3.78, 4.41, 4.51, 4.51, 4.51, 4.51, 3.89, 1.49, 0.06, 0.04)
    7  = @(50,  3.78, 4.41, 4.5,  4.51, 4.51, 3.79, 3.26, 1.5,  0.05, 0.05)
    8  = @(40,  3.78, 4.41, 4.51, 4.51, 3.64, 3.05, 2.63, 1.51, 0.06, 0.06)
    9  = @(30,  3.78, 4.41, 4.51, 3.42, 2.75, 2.31, 1.98, 1.8,  1.8,  1.8)
    10 = @(20,  3.78, 4.42, 3.05, 2.31, 1.86, 1.8,  1.8,  1.8,  1.81, 1.81)
    11 = @(10,  3.77, 2.28, 1.77, 1.77, 1.77, 1.77, 1.77, 1.77, 1.77, 1.77)
}

$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# New blank, number-formatted placeholder cells (columns M:V for rows
# 2-11, with row 2 extending one extra column to W) reserved for the
# upcoming nodule/AMF allocation figure data. Applying a "0.00" number
# format creates the new cellXfs entry (numFmtId 2) used by these cells.
# ---------------------------------------------------------------------
$ws.Range("M2:W2").NumberFormat = "0.00"
for ($r = 3; $r -le 11; $r++) {
    $ws.Range("M$($r):V$($r)").NumberFormat = "0.00"
}

# ---------------------------------------------------------------------
# Update the current selection to match where the author left off
# (selecting the new placeholder rows 13-22 across all columns).
# ---------------------------------------------------------------------
$ws.Range("A13:XFD22").Select()
